$d = $word.ActiveDocument

# --- 1) Split the run containing "...9" into two runs: "...(1..." | "9" ---
# This mirrors the author's edit where the trailing "9" of
# "чество заданных номеров (1...9" ends up as its own run (visible text is
# unchanged: "чество заданных номеров (1...9").
$rng = $d.Content
$found = $rng.Find.Execute("номеров (1...9", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$nineRange = $d.Range($rng.End - 1, $rng.End)
$nineRange.Font.Bold = $true
$nineRange.Font.Bold = $false

# --- 2) Append a new "Меню" Heading 1 paragraph, followed by an empty
#        paragraph, at the very end of the document body. ---
$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.Text = "Меню"

$menuPara = $d.Paragraphs.Last
$menuPara.Style = "Heading 1"

$endRange2 = $d.Range($d.Content.End, $d.Content.End)
$endRange2.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>") | Out-Null
